$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---- Update title/header shared strings (Volume/Number and report week dates) ----
$ws.Cells.Replace("Number  22", "Number  23") | Out-Null
$ws.Cells.Replace("5/29/2023", "6/5/2023") | Out-Null
$ws.Cells.Replace("6/4/2023", "6/11/2023") | Out-Null

# ---- Universal donor cells (their own style/value never change) ----
# C14 = style14 text '0' (shared string idx 20)
# M14 = style14 text '***.*' (shared string idx 21)
# I14 = style15 numeric (#,##0)
# L14 = style16 numeric percent (#,##0.0)
$donor20 = $ws.Cells.Item(14, 3)
$donor21 = $ws.Cells.Item(14, 13)
$donor15 = $ws.Cells.Item(14, 9)
$donor16 = $ws.Cells.Item(14, 12)

# ---- Cells converting from text (dash-style) to numeric: fix style via donor copy, then set value ----
$donor15.Copy($ws.Cells.Item(14, 4))
$ws.Cells.Item(14, 4).Value = 2
$donor16.Copy($ws.Cells.Item(14, 5))
$ws.Cells.Item(14, 5).Value = -100
$donor15.Copy($ws.Cells.Item(14, 7))
$ws.Cells.Item(14, 7).Value = 2
$donor16.Copy($ws.Cells.Item(14, 8))
$ws.Cells.Item(14, 8).Value = -100
$donor15.Copy($ws.Cells.Item(22, 4))
$ws.Cells.Item(22, 4).Value = 1
$donor16.Copy($ws.Cells.Item(22, 5))
$ws.Cells.Item(22, 5).Value = -100
$donor15.Copy($ws.Cells.Item(22, 7))
$ws.Cells.Item(22, 7).Value = 1
$donor16.Copy($ws.Cells.Item(22, 8))
$ws.Cells.Item(22, 8).Value = -100
$donor15.Copy($ws.Cells.Item(26, 4))
$ws.Cells.Item(26, 4).Value = 3
$donor16.Copy($ws.Cells.Item(26, 5))
$ws.Cells.Item(26, 5).Value = -33.333333333333

# ---- Cells converting from numeric to text (dash-style): fix via donor copy (value already correct) ----
$donor20.Copy($ws.Cells.Item(15, 3))
$donor20.Copy($ws.Cells.Item(22, 6))
$donor20.Copy($ws.Cells.Item(23, 7))
$donor21.Copy($ws.Cells.Item(23, 8))
$donor20.Copy($ws.Cells.Item(27, 3))
$donor20.Copy($ws.Cells.Item(28, 6))
$donor20.Copy($ws.Cells.Item(29, 6))

# ---- Simple numeric value updates (style unchanged) ----
$ws.Cells.Item(14, 10).Value = 9
$ws.Cells.Item(14, 11).Value = 0
$ws.Cells.Item(14, 14).Value = -50
$ws.Cells.Item(15, 6).Value = 2
$ws.Cells.Item(15, 7).Value = 1
$ws.Cells.Item(15, 8).Value = 100
$ws.Cells.Item(15, 9).Value = 14
$ws.Cells.Item(15, 11).Value = 0
$ws.Cells.Item(15, 12).Value = 0
$ws.Cells.Item(15, 13).Value = 40
$ws.Cells.Item(15, 14).Value = -53.333333333333
$ws.Cells.Item(16, 3).Value = 5
$ws.Cells.Item(16, 4).Value = 10
$ws.Cells.Item(16, 5).Value = -50
$ws.Cells.Item(16, 6).Value = 31
$ws.Cells.Item(16, 7).Value = 38
$ws.Cells.Item(16, 8).Value = -18.421052631578
$ws.Cells.Item(16, 9).Value = 155
$ws.Cells.Item(16, 10).Value = 190
$ws.Cells.Item(16, 11).Value = -18.421052631578
$ws.Cells.Item(16, 12).Value = 20.155038759689
$ws.Cells.Item(16, 13).Value = 1.307189542483
$ws.Cells.Item(16, 14).Value = -69.367588932806
$ws.Cells.Item(17, 3).Value = 7
$ws.Cells.Item(17, 4).Value = 22
$ws.Cells.Item(17, 5).Value = -68.181818181818
$ws.Cells.Item(17, 6).Value = 38
$ws.Cells.Item(17, 7).Value = 57
$ws.Cells.Item(17, 8).Value = -33.333333333333
$ws.Cells.Item(17, 9).Value = 281
$ws.Cells.Item(17, 10).Value = 263
$ws.Cells.Item(17, 11).Value = 6.844106463878
$ws.Cells.Item(17, 12).Value = 44.845360824742
$ws.Cells.Item(17, 13).Value = 91.156462585034
$ws.Cells.Item(17, 14).Value = -3.103448275862
$ws.Cells.Item(18, 3).Value = 6
$ws.Cells.Item(18, 5).Value = 50
$ws.Cells.Item(18, 6).Value = 24
$ws.Cells.Item(18, 7).Value = 20
$ws.Cells.Item(18, 8).Value = 20
$ws.Cells.Item(18, 9).Value = 97
$ws.Cells.Item(18, 10).Value = 130
$ws.Cells.Item(18, 11).Value = -25.384615384615
$ws.Cells.Item(18, 12).Value = 22.784810126582
$ws.Cells.Item(18, 13).Value = -3
$ws.Cells.Item(18, 14).Value = -83.503401360544
$ws.Cells.Item(19, 3).Value = 12
$ws.Cells.Item(19, 4).Value = 8
$ws.Cells.Item(19, 5).Value = 50
$ws.Cells.Item(19, 6).Value = 46
$ws.Cells.Item(19, 7).Value = 28
$ws.Cells.Item(19, 8).Value = 64.285714285714
$ws.Cells.Item(19, 9).Value = 210
$ws.Cells.Item(19, 10).Value = 200
$ws.Cells.Item(19, 11).Value = 5
$ws.Cells.Item(19, 12).Value = 16.022099447513
$ws.Cells.Item(19, 13).Value = 62.790697674418
$ws.Cells.Item(19, 14).Value = 7.142857142857
$ws.Cells.Item(20, 9).Value = 154
$ws.Cells.Item(20, 10).Value = 123
$ws.Cells.Item(20, 11).Value = 25.203252032520
$ws.Cells.Item(20, 12).Value = 126.470588235294
$ws.Cells.Item(20, 13).Value = 250
$ws.Cells.Item(20, 14).Value = -41.666666666666
$ws.Cells.Item(21, 3).Value = 36
$ws.Cells.Item(21, 4).Value = 52
$ws.Cells.Item(21, 5).Value = -30.769230769230
$ws.Cells.Item(21, 6).Value = 166
$ws.Cells.Item(21, 7).Value = 174
$ws.Cells.Item(21, 8).Value = -4.597701149425
$ws.Cells.Item(21, 9).Value = 920
$ws.Cells.Item(21, 10).Value = 929
$ws.Cells.Item(21, 11).Value = -0.968783638320
$ws.Cells.Item(21, 12).Value = 36.904761904761
$ws.Cells.Item(21, 13).Value = 57.804459691252
$ws.Cells.Item(21, 14).Value = -51.374207188160
$ws.Cells.Item(22, 10).Value = 3
$ws.Cells.Item(22, 11).Value = -33.333333333333
$ws.Cells.Item(22, 12).Value = -60
$ws.Cells.Item(24, 4).Value = 23
$ws.Cells.Item(24, 5).Value = 0
$ws.Cells.Item(24, 6).Value = 88
$ws.Cells.Item(24, 7).Value = 107
$ws.Cells.Item(24, 8).Value = -17.757009345794
$ws.Cells.Item(24, 9).Value = 467
$ws.Cells.Item(24, 10).Value = 497
$ws.Cells.Item(24, 11).Value = -6.036217303822
$ws.Cells.Item(24, 12).Value = 27.595628415300
$ws.Cells.Item(24, 13).Value = 17.632241813602
$ws.Cells.Item(25, 3).Value = 12
$ws.Cells.Item(25, 4).Value = 12
$ws.Cells.Item(25, 5).Value = 0
$ws.Cells.Item(25, 6).Value = 88
$ws.Cells.Item(25, 7).Value = 61
$ws.Cells.Item(25, 8).Value = 44.262295081967
$ws.Cells.Item(25, 9).Value = 375
$ws.Cells.Item(25, 10).Value = 401
$ws.Cells.Item(25, 11).Value = -6.483790523690
$ws.Cells.Item(25, 12).Value = 25.838926174496
$ws.Cells.Item(25, 13).Value = -11.971830985915
$ws.Cells.Item(26, 7).Value = 5
$ws.Cells.Item(26, 8).Value = 20
$ws.Cells.Item(26, 9).Value = 26
$ws.Cells.Item(26, 10).Value = 25
$ws.Cells.Item(26, 11).Value = 4
$ws.Cells.Item(26, 12).Value = 18.181818181818
$ws.Cells.Item(27, 4).Value = 2
$ws.Cells.Item(27, 5).Value = -100
$ws.Cells.Item(27, 6).Value = 3
$ws.Cells.Item(27, 7).Value = 6
$ws.Cells.Item(27, 8).Value = -50
$ws.Cells.Item(27, 9).Value = 29
$ws.Cells.Item(27, 10).Value = 35
$ws.Cells.Item(27, 11).Value = -17.142857142857
$ws.Cells.Item(27, 12).Value = 16
$ws.Cells.Item(28, 4).Value = 3
$ws.Cells.Item(28, 7).Value = 6
$ws.Cells.Item(28, 8).Value = -100
$ws.Cells.Item(28, 10).Value = 26
$ws.Cells.Item(28, 11).Value = -30.769230769230
$ws.Cells.Item(28, 12).Value = -33.333333333333
$ws.Cells.Item(28, 13).Value = 12.5
$ws.Cells.Item(29, 4).Value = 3
$ws.Cells.Item(29, 7).Value = 6
$ws.Cells.Item(29, 8).Value = -100
$ws.Cells.Item(29, 10).Value = 23
$ws.Cells.Item(29, 11).Value = -52.173913043478
$ws.Cells.Item(29, 12).Value = -57.692307692307
$ws.Cells.Item(29, 13).Value = -15.384615384615

Write-Output "edit complete"